$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3116")

# Row 10: Objetivos -> teacher name (Carlos)
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# Row 13: Programa resumido -> activation date (must stay plain text, not be
# auto-converted to a date serial). Writing it as a formula string and then
# pasting-special as values keeps the original (General) cell style while
# still landing as a text value, exactly like the other label/value pairs.
$ws.Range("B13").Formula = '="01/01/2023"'
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4163) | Out-Null
$ws.Range("C13").Formula = '="01/01/2023"'
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4163) | Out-Null

# Row 15: Programa -> teacher name (Carlos) again
$ws.Range("B15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# Row 18: Método -> teacher name (Cassius)
$ws.Range("B18").Value = "3586455 - Cassius Olivio Figueiredo Terra Ruchert"
$ws.Range("C18").Value = "3586455 - Cassius Olivio Figueiredo Terra Ruchert"

$excel.CutCopyMode = 0
